$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '38.255.68'
$ws.Range("E2").Value = '  +2.80%  '
$ws.Range("D3").Value = '2.056.93'
$ws.Range("E3").Value = '  +1.70%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("D7").Value = "'61.02"
$ws.Range("E7").Value = '  +8.89%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +2.37%  '
$ws.Range("D10").Value = "'0.0827"
$ws.Range("E10").Value = '  +5.67%  '
$ws.Range("E11").Value = '  +1.92%  '
$ws.Range("D12").Value = "'14.85"
$ws.Range("E12").Value = '  +3.99%  '
$ws.Range("D13").Value = '2.361.47'
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("D14").Value = "'21.16"
$ws.Range("E14").Value = '  +4.81%  '
$ws.Range("D15").Value = "'0.764"
$ws.Range("E15").Value = '  +3.33%  '
$ws.Range("E16").Value = '  +2.89%  '
$ws.Range("D17").Value = '2.063.94'
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = '38.208.12'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("E19").Value = '  +1.63%  '
$ws.Range("D20").Value = "'69.88"
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("D21").Value = '0.0₃0836'
$ws.Range("E21").Value = '  +2.32%  '
$ws.Range("D22").Value = "'225.38"
$ws.Range("E22").Value = '  +1.03%  '
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -0.57%  '
$ws.Range("D25").Value = "'2.22"
$ws.Range("E25").Value = '  -0.76%  '
$ws.Range("D26").Value = "'9.29"
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").Value = "'166.36"
$ws.Range("E27").Value = '  +1.71%  '
$ws.Range("D28").Value = "'0.132"
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = "'18.99"
$ws.Range("E29").Value = '  +1.57%  '
$ws.Range("D30").Value = "'1.29"
$ws.Range("E30").Value = '  -1.32%  '
$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = '  +2.56%  '
$ws.Range("D32").Value = "'4.51"
$ws.Range("E32").Value = '  +1.30%  '
$ws.Range("D33").Value = "'4.58"
$ws.Range("E33").Value = '  +3.21%  '
$ws.Range("E34").Value = '  +3.28%  '
$ws.Range("D35").Value = "'0.0602"
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("D36").Value = "'6.39"
$ws.Range("E36").Value = '  +14.84%  '
$ws.Range("E37").Value = '  -1.39%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = '1.529.93'
$ws.Range("E40").Value = '  +4.26%  '
$ws.Range("D41").Value = "'97.85"
$ws.Range("E41").Value = '  +4.17%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = "'0.0217"
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("B43").Value = 'InjectiveProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D43").Value = "'16.93"
$ws.Range("E43").Value = '  +4.44%  '
$ws.Range("E44").Value = '  +1.43%  '
$ws.Range("E45").Value = '  +1.96%  '
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("D47").Value = "'4.01"
$ws.Range("E47").Value = '  -8.78%  '
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("D49").Value = "'2.98"
$ws.Range("E49").Value = '  +1.93%  '
$ws.Range("D50").Value = "'7.05"
$ws.Range("E50").Value = '  -0.93%  '
$ws.Range("D51").Value = '2.251.10'
$ws.Range("E51").Value = '  +1.76%  '
